$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-12-09 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-10 Wednesday", 2)

# Update the multiplication table cells by explicit (row, column) address
# so that duplicate "old" values (e.g. 480×6=2880) are replaced independently.
$tbl = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "483×3=1449" },
    @{ Row = 1;  Col = 2; Text = "213×9=1917" },
    @{ Row = 1;  Col = 3; Text = "404×9=3636" },
    @{ Row = 1;  Col = 4; Text = "232×9=2088" },
    @{ Row = 1;  Col = 5; Text = "292×4=1168" },

    @{ Row = 5;  Col = 1; Text = "145×7=1015" },
    @{ Row = 5;  Col = 2; Text = "252×6=1512" },
    @{ Row = 5;  Col = 3; Text = "147×6=882" },
    @{ Row = 5;  Col = 4; Text = "726×8=5808" },
    @{ Row = 5;  Col = 5; Text = "119×9=1071" },

    @{ Row = 10; Col = 1; Text = "533×5=2665" },
    @{ Row = 10; Col = 2; Text = "503×7=3521" },
    @{ Row = 10; Col = 3; Text = "933×7=6531" },
    @{ Row = 10; Col = 4; Text = "757×7=5299" },
    @{ Row = 10; Col = 5; Text = "901×3=2703" },

    @{ Row = 15; Col = 1; Text = "448×8=3584" },
    @{ Row = 15; Col = 2; Text = "521×5=2605" },
    @{ Row = 15; Col = 3; Text = "978×7=6846" },
    @{ Row = 15; Col = 4; Text = "616×3=1848" },
    @{ Row = 15; Col = 5; Text = "249×8=1992" },

    @{ Row = 20; Col = 1; Text = "202×8=1616" },
    @{ Row = 20; Col = 2; Text = "678×3=2034" },
    @{ Row = 20; Col = 3; Text = "118×3=354" },
    @{ Row = 20; Col = 4; Text = "282×4=1128" },
    @{ Row = 20; Col = 5; Text = "245×6=1470" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
